$d = $word.ActiveDocument

# Locate the paragraph that ends the bibliography entry
# ("Janeiro: Editora Interciência , 2004."). The three paragraphs that
# follow it -- a blank paragraph, the "Ver no Jupiter..." paragraph, and
# the "© 2020 ..." footer paragraph -- are removed by this edit, while the
# blank paragraph and page-break paragraph that come after them are kept.

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Janeiro: Editora Interci*2004.*") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Anchor paragraph not found"
}

$blank1 = $anchor.Next()
$jupiter = $blank1.Next()
$copyright = $jupiter.Next()

if ($jupiter.Range.Text -notlike "*Ver no Jupiter*") {
    throw "Unexpected paragraph where 'Ver no Jupiter' text was expected: $($jupiter.Range.Text)"
}
if ($copyright.Range.Text -notlike "*2020*Creative Commons*") {
    throw "Unexpected paragraph where copyright text was expected: $($copyright.Range.Text)"
}

$start = $blank1.Range.Start
$end = $copyright.Range.End

$r = $d.Range($start, $end)
$r.Delete()
